# Apply crypto price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells that receive numeric-looking text to stay as Text,
# matching the original inlineStr/text cell type in the workbook.
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D48").NumberFormat = "@"
$ws.Range("D50:D51").NumberFormat = "@"

$ws.Range("D2").Value = '62.377.26'
$ws.Range("E2").Value = '  +2.21%  '

$ws.Range("D3").Value = '2.922.16'
$ws.Range("E3").Value = '  +1.37%  '

$ws.Range("E4").Value = '  +0.30%  '

$ws.Range("D5").Value = '588.27'
$ws.Range("E5").Value = '  -0.21%  '

$ws.Range("D6").Value = '146.86'
$ws.Range("E6").Value = '  +5.43%  '

$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("D8").Value = '0.504'
$ws.Range("E8").Value = '  +2.26%  '

$ws.Range("D9").Value = '2.912.19'
$ws.Range("E9").Value = '  +1.01%  '

$ws.Range("D10").Value = '7.07'
$ws.Range("E10").Value = '  -0.77%  '

$ws.Range("D11").Value = '0.150'
$ws.Range("E11").Value = '  +8.67%  '

$ws.Range("D12").Value = '0.433'
$ws.Range("E12").Value = '  +1.16%  '

$ws.Range("D13").Value = '0.0000235'
$ws.Range("E13").Value = '  +8.01%  '

$ws.Range("D14").Value = '32.11'
$ws.Range("E14").Value = '  -0.42%  '

$ws.Range("D15").Value = '0.125'
$ws.Range("E15").Value = '  -0.69%  '

$ws.Range("D16").Value = '3.423.36'
$ws.Range("E16").Value = '  +1.84%  '

$ws.Range("D17").Value = '62.512.96'
$ws.Range("E17").Value = '  +2.65%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.930.17'
$ws.Range("E18").Value = '  +1.98%  '

$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '6.58'
$ws.Range("E19").Value = '  +1.14%  '

$ws.Range("D20").Value = '433.80'
$ws.Range("E20").Value = '  +1.73%  '

$ws.Range("D21").Value = '13.33'
$ws.Range("E21").Value = '  +1.26%  '

$ws.Range("D22").Value = '0.657'
$ws.Range("E22").Value = '  +0.35%  '

$ws.Range("D23").Value = '6.88'
$ws.Range("E23").Value = '  -0.31%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '80.06'
$ws.Range("E24").Value = '  +0.19%  '

$ws.Range("B25").Value = 'RenderToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D25").Value = '11.01'
$ws.Range("E25").Value = '  +5.17%  '

$ws.Range("D26").Value = '11.77'
$ws.Range("E26").Value = '  +3.12%  '

$ws.Range("D27").Value = '2.09'
$ws.Range("E27").Value = '  +0.82%  '

$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.17%  '

$ws.Range("D29").Value = '7.21'
$ws.Range("E29").Value = '  +8.95%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '2.58'
$ws.Range("E30").Value = '  +1.50%  '

$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Value = '0.0000102'
$ws.Range("E31").Value = '  +21.08%  '

$ws.Range("D32").Value = '2.12'
$ws.Range("E32").Value = '  +2.31%  '

$ws.Range("D33").Value = '0.108'
$ws.Range("E33").Value = '  +4.10%  '

$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.28%  '

$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = '25.93'
$ws.Range("E35").Value = '  +1.15%  '

$ws.Range("D36").Value = '0.975'
$ws.Range("E36").Value = '  +0.87%  '

$ws.Range("D37").Value = '3.07'
$ws.Range("E37").Value = '  +9.91%  '

$ws.Range("D38").Value = '5.49'
$ws.Range("E38").Value = '  +1.29%  '

$ws.Range("D39").Value = '49.38'
$ws.Range("E39").Value = '  +0.93%  '

$ws.Range("D40").Value = '2.00'
$ws.Range("E40").Value = '  +5.11%  '

$ws.Range("D41").Value = '8.33'
$ws.Range("E41").Value = '  +0.03%  '

$ws.Range("D42").Value = '0.115'
$ws.Range("E42").Value = '  -0.19%  '

$ws.Range("D43").Value = '0.275'
$ws.Range("E43").Value = '  +3.22%  '

$ws.Range("D44").Value = '39.48'
$ws.Range("E44").Value = '  +2.31%  '

$ws.Range("D45").Value = '136.05'
$ws.Range("E45").Value = '  +3.49%  '

$ws.Range("D46").Value = '2.690.16'
$ws.Range("E46").Value = '  +1.13%  '

$ws.Range("D47").Value = '0.0338'
$ws.Range("E47").Value = '  +2.61%  '

$ws.Range("D48").Value = '353.24'
$ws.Range("E48").Value = '  -0.74%  '

$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("D50").Value = '0.103'
$ws.Range("E50").Value = '  +1.71%  '

$ws.Range("D51").Value = '22.40'
$ws.Range("E51").Value = '  +0.64%  '
